$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "wefwef"
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Range("A4").Select()
